$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new blank cells D1,E1,D2,E2 (default/"Normal" style, same as existing C4:C11 style)
$ws.Range("C4").Copy()
$ws.Range("D1:E2").PasteSpecial(-4122)

# Row 3 ("Combo box Passengers"): extend with additional passenger count values
# matching the format/style already used by B3 (text number format)
$ws.Range("B3").Copy()
$ws.Range("C3:E3").PasteSpecial(-4122)

$ws.Range("C3").Value = "2"
$ws.Range("D3").Value = "3"
$ws.Range("E3").Value = "4"

# Column widths: C narrower, D and E new columns
# (ColumnWidth values chosen so that, after Excel's internal pixel
# rounding, the serialized sheet width is as close as possible to the
# target widths 6.14814814814815 / 5.03703703703704 / 6.85185185185185)
$ws.Columns.Item(3).ColumnWidth = 5.333333333333333
$ws.Columns.Item(4).ColumnWidth = 4.166666666666667
$ws.Columns.Item(5).ColumnWidth = 6

$excel.CutCopyMode = $false

# Selection moves to F3
$ws.Range("F3").Select()
